# Update countries & provincias Spain
# Refresh the COVID stats table on sheet "Pais": update the case counters
# for the countries whose figures changed, and re-seat the country names
# for the four rank-swap pairs produced by the refreshed sort order
# (Suecia/Ucrania, Kirguistan/Armenia, Republica de Chipre/Georgia,
# Islas Malvinas/Montserrat). Also bump the "Datos actualizados" footer
# timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: Estados Unidos (no reorder, just new totals) ---
$ws.Cells.Item(4, 4).Value = 2755608
$ws.Cells.Item(4, 5).Value = 2382600

# --- Row 6: India (no reorder, just new totals) ---
$ws.Cells.Item(6, 2).Value = 2332908
$ws.Cells.Item(6, 3).Value = 4503
$ws.Cells.Item(6, 4).Value = 1640362
$ws.Cells.Item(6, 5).Value = 646330
$ws.Cells.Item(6, 7).Value = 28
$ws.Cells.Item(6, 8).Value = 46216

# --- Rows 35/36: Ucrania overtakes Suecia ---
$ws.Cells.Item(35, 1).Value = "Ucrania"
$ws.Cells.Item(35, 2).Value = 84548
$ws.Cells.Item(35, 3).Value = 1433
$ws.Cells.Item(35, 4).Value = 45686
$ws.Cells.Item(35, 5).Value = 36892
$ws.Cells.Item(35, 7).Value = 19
$ws.Cells.Item(35, 8).Value = 1970

$ws.Cells.Item(36, 1).Value = "Suecia"
$ws.Cells.Item(36, 2).Value = 83126
$ws.Cells.Item(36, 4).Value = 0
$ws.Cells.Item(36, 5).Value = 0
$ws.Cells.Item(36, 8).Value = 5770

# --- Rows 55/56: Armenia overtakes Kirguistan ---
$ws.Cells.Item(55, 1).Value = "Armenia"
$ws.Cells.Item(55, 2).Value = 40794
$ws.Cells.Item(55, 3).Value = 201
$ws.Cells.Item(55, 4).Value = 33492
$ws.Cells.Item(55, 5).Value = 6496
$ws.Cells.Item(55, 7).Value = 3
$ws.Cells.Item(55, 8).Value = 806

$ws.Cells.Item(56, 1).Value = "Kirguistan"
$ws.Cells.Item(56, 2).Value = 40759
$ws.Cells.Item(56, 3).Value = 304
$ws.Cells.Item(56, 4).Value = 32997
$ws.Cells.Item(56, 5).Value = 6278
$ws.Cells.Item(56, 7).Value = 6
$ws.Cells.Item(56, 8).Value = 1484

# --- Row 57: Afganistan (no reorder, just new totals) ---
$ws.Cells.Item(57, 2).Value = 37345
$ws.Cells.Item(57, 3).Value = 76
$ws.Cells.Item(57, 4).Value = 26694
$ws.Cells.Item(57, 5).Value = 9297
$ws.Cells.Item(57, 7).Value = 10
$ws.Cells.Item(57, 8).Value = 1354

# --- Row 58: Suiza (no reorder, just new totals) ---
$ws.Cells.Item(58, 4).Value = 32700
$ws.Cells.Item(58, 5).Value = 2205

# --- Row 107: Zimbabue (no reorder, just new totals) ---
$ws.Cells.Item(107, 4).Value = 1544
$ws.Cells.Item(107, 5).Value = 3170

# --- Row 108: Hungria (no reorder, just new totals) ---
$ws.Cells.Item(108, 2).Value = 4768
$ws.Cells.Item(108, 3).Value = 22
$ws.Cells.Item(108, 4).Value = 3529
$ws.Cells.Item(108, 5).Value = 634

# --- Row 122: Sri Lanka (no reorder, just new totals) ---
$ws.Cells.Item(122, 4).Value = 2638
$ws.Cells.Item(122, 5).Value = 231

# --- Row 136: Yemen (no reorder, just new totals) ---
$ws.Cells.Item(136, 4).Value = 919
$ws.Cells.Item(136, 5).Value = 389

# --- Row 144: Letonia (no reorder, just new totals) ---
$ws.Cells.Item(144, 2).Value = 1303
$ws.Cells.Item(144, 3).Value = 10
$ws.Cells.Item(144, 5).Value = 193

# --- Rows 146/147: Georgia overtakes Republica de Chipre ---
$ws.Cells.Item(146, 1).Value = "Georgia"
$ws.Cells.Item(146, 2).Value = 1278
$ws.Cells.Item(146, 3).Value = 14
$ws.Cells.Item(146, 4).Value = 1058
$ws.Cells.Item(146, 5).Value = 203
$ws.Cells.Item(146, 8).Value = 17

$ws.Cells.Item(147, 1).Value = "Republica de Chipre"
$ws.Cells.Item(147, 2).Value = 1277
$ws.Cells.Item(147, 4).Value = 870
$ws.Cells.Item(147, 5).Value = 387
$ws.Cells.Item(147, 8).Value = 20

# --- Row 167: Taiwan (no reorder, just new totals) ---
$ws.Cells.Item(167, 2).Value = 481
$ws.Cells.Item(167, 3).Value = 1
$ws.Cells.Item(167, 4).Value = 450
$ws.Cells.Item(167, 5).Value = 24

# --- Row 186: Bermudas (no reorder, just new totals) ---
$ws.Cells.Item(186, 4).Value = 145
$ws.Cells.Item(186, 5).Value = 5

# --- Rows 213/214: Montserrat overtakes Islas Malvinas ---
$ws.Cells.Item(213, 1).Value = "Montserrat"
$ws.Cells.Item(213, 4).Value = 12
$ws.Cells.Item(213, 8).Value = 1

$ws.Cells.Item(214, 1).Value = "Islas Malvinas"
$ws.Cells.Item(214, 4).Value = 13
$ws.Cells.Item(214, 8).Value = 0

# --- Footer: refresh the "last updated" timestamp ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 12 de Agosto de 2020 a las 09:29"
